$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update each data row (hours 1-17): date 2024-08-14 -> 2024-08-15, plus refreshed load figures
# Row 2 (Hour 1)
$ws.Range("A2").Value = "'2024-08-15"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 22266
$ws.Range("D2").Value = 5999
$ws.Range("E2").Value = 9984
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 6030
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 17786
$ws.Range("J2").Value = 16906
$ws.Range("K2").Value = 78971
$ws.Range("L2").Value = 22268.2266
$ws.Range("M2").Value = 5999.5999
$ws.Range("N2").Value = 9984.9984
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 6030.603
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 17787.7786
$ws.Range("S2").Value = 16907.6906
$ws.Range("T2").Value = 78978.8971

# Row 3 (Hour 2)
$ws.Range("A3").Value = "'2024-08-15"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 21326
$ws.Range("D3").Value = 5717
$ws.Range("E3").Value = 9532
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 5923
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 17126
$ws.Range("J3").Value = 16079
$ws.Range("K3").Value = 75703
$ws.Range("L3").Value = 21328.1326
$ws.Range("M3").Value = 5717.5717
$ws.Range("N3").Value = 9532.9532
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 5923.5923
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 17127.7126
$ws.Range("S3").Value = 16080.6079
$ws.Range("T3").Value = 75710.5703

# Row 4 (Hour 3)
$ws.Range("A4").Value = "'2024-08-15"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 20449
$ws.Range("D4").Value = 5454
$ws.Range("E4").Value = 9230
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 5604
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 16410
$ws.Range("J4").Value = 15490
$ws.Range("K4").Value = 72637
$ws.Range("L4").Value = 20451.0449
$ws.Range("M4").Value = 5454.5454
$ws.Range("N4").Value = 9230.923
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 5604.5604
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 16411.641
$ws.Range("S4").Value = 15491.549
$ws.Range("T4").Value = 72644.2637

# Row 5 (Hour 4)
$ws.Range("A5").Value = "'2024-08-15"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 19762
$ws.Range("D5").Value = 5204
$ws.Range("E5").Value = 8913
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 5557
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 15487
$ws.Range("J5").Value = 15148
$ws.Range("K5").Value = 70071
$ws.Range("L5").Value = 19763.9762
$ws.Range("M5").Value = 5204.5204
$ws.Range("N5").Value = 8913.8913
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 5557.5557
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 15488.5487
$ws.Range("S5").Value = 15149.5148
$ws.Range("T5").Value = 70078.0071

# Row 6 (Hour 5)
$ws.Range("A6").Value = "'2024-08-15"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 19910
$ws.Range("D6").Value = 5266
$ws.Range("E6").Value = 8941
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 5670
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 15633
$ws.Range("J6").Value = 15424
$ws.Range("K6").Value = 70844
$ws.Range("L6").Value = 19911.991
$ws.Range("M6").Value = 5266.5266
$ws.Range("N6").Value = 8941.8941
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 5670.567
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 15634.5633
$ws.Range("S6").Value = 15425.5424
$ws.Range("T6").Value = 70851.0844

# Row 7 (Hour 6)
$ws.Range("A7").Value = "'2024-08-15"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 19707
$ws.Range("D7").Value = 5468
$ws.Range("E7").Value = 9267
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 7157
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 16091
$ws.Range("J7").Value = 15868
$ws.Range("K7").Value = 73558
$ws.Range("L7").Value = 19708.9707
$ws.Range("M7").Value = 5468.5468
$ws.Range("N7").Value = 9267.9267
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 7157.7157
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 16092.6091
$ws.Range("S7").Value = 15869.5868
$ws.Range("T7").Value = 73565.3558

# Row 8 (Hour 7)
$ws.Range("A8").Value = "'2024-08-15"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 20268
$ws.Range("D8").Value = 4604
$ws.Range("E8").Value = 9666
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 7401
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 15557
$ws.Range("J8").Value = 16013
$ws.Range("K8").Value = 73509
$ws.Range("L8").Value = 20270.0268
$ws.Range("M8").Value = 4604.4604
$ws.Range("N8").Value = 9666.9666
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 7401.7401
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 15558.5557
$ws.Range("S8").Value = 16014.6013
$ws.Range("T8").Value = 73516.3509

# Row 9 (Hour 8)
$ws.Range("A9").Value = "'2024-08-15"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 23185
$ws.Range("D9").Value = 4901
$ws.Range("E9").Value = 11223
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 10623
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 18375
$ws.Range("J9").Value = 18064
$ws.Range("K9").Value = 86371
$ws.Range("L9").Value = 23187.3185
$ws.Range("M9").Value = 4901.4901
$ws.Range("N9").Value = 11224.1223
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 10624.0623
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 18376.8375
$ws.Range("S9").Value = 18065.8064
$ws.Range("T9").Value = 86379.6371

# Row 10 (Hour 9)
$ws.Range("A10").Value = "'2024-08-15"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 28823
$ws.Range("D10").Value = 5818
$ws.Range("E10").Value = 13930
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 14092
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 20118
$ws.Range("J10").Value = 21041
$ws.Range("K10").Value = 103822
$ws.Range("L10").Value = 28825.8823
$ws.Range("M10").Value = 5818.5818
$ws.Range("N10").Value = 13931.393
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 14093.4092
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 20120.0118
$ws.Range("S10").Value = 21043.1041
$ws.Range("T10").Value = 103832.3822

# Row 11 (Hour 10)
$ws.Range("A11").Value = "'2024-08-15"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 32025
$ws.Range("D11").Value = 6084
$ws.Range("E11").Value = 18165
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 15134
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 22980
$ws.Range("K11").Value = 94388
$ws.Range("L11").Value = 32028.2025
$ws.Range("M11").Value = 6084.6084
$ws.Range("N11").Value = 18166.8165
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 15135.5134
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 22982.298
$ws.Range("T11").Value = 94397.4388

# Row 12 (Hour 11)
$ws.Range("A12").Value = "'2024-08-15"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 32138
$ws.Range("D12").Value = 6317
$ws.Range("E12").Value = 19577
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 15265
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 23801
$ws.Range("K12").Value = 97098
$ws.Range("L12").Value = 32141.2138
$ws.Range("M12").Value = 6317.6317
$ws.Range("N12").Value = 19578.9577
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 15266.5265
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 0
$ws.Range("S12").Value = 23803.3801
$ws.Range("T12").Value = 97107.7098

# Row 13 (Hour 12)
$ws.Range("A13").Value = "'2024-08-15"
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 31865
$ws.Range("D13").Value = 6389
$ws.Range("E13").Value = 19812
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 15463
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 23781
$ws.Range("K13").Value = 97310
$ws.Range("L13").Value = 31868.1865
$ws.Range("M13").Value = 6389.6389
$ws.Range("N13").Value = 19813.9812
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 15464.5463
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 23783.3781
$ws.Range("T13").Value = 97319.731

# Row 14 (Hour 13)
$ws.Range("A14").Value = "'2024-08-15"
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").Value = 13
$ws.Range("C14").Value = 32062
$ws.Range("D14").Value = 6461
$ws.Range("E14").Value = 19870
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 15114
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 23508
$ws.Range("K14").Value = 97015
$ws.Range("L14").Value = 32065.2062
$ws.Range("M14").Value = 6461.6461
$ws.Range("N14").Value = 19871.987
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 15115.5114
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 23510.3508
$ws.Range("T14").Value = 97024.7015

# Row 15 (Hour 14)
$ws.Range("A15").Value = "'2024-08-15"
$ws.Range("A15").Style = "Normal"
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = 33756
$ws.Range("D15").Value = 6996
$ws.Range("E15").Value = 20143
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 16097
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 22648
$ws.Range("J15").Value = 24132
$ws.Range("K15").Value = 123772
$ws.Range("L15").Value = 33759.3756
$ws.Range("M15").Value = 6996.6996
$ws.Range("N15").Value = 20145.0143
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 16098.6097
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 22650.2648
$ws.Range("S15").Value = 24134.4132
$ws.Range("T15").Value = 123784.3772

# Row 16 (Hour 15)
$ws.Range("A16").Value = "'2024-08-15"
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = 33340
$ws.Range("D16").Value = 6897
$ws.Range("E16").Value = 19871
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 15846
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 22412
$ws.Range("J16").Value = 24607
$ws.Range("K16").Value = 122973
$ws.Range("L16").Value = 33343.334
$ws.Range("M16").Value = 6897.6897
$ws.Range("N16").Value = 19872.9871
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 15847.5846
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 22414.2412
$ws.Range("S16").Value = 24609.4607
$ws.Range("T16").Value = 122985.2973

# Row 17 (Hour 16)
$ws.Range("A17").Value = "'2024-08-15"
$ws.Range("A17").Style = "Normal"
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = 33211
$ws.Range("D17").Value = 6639
$ws.Range("E17").Value = 20137
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 22106
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 82093
$ws.Range("L17").Value = 33214.3211
$ws.Range("M17").Value = 6639.6639
$ws.Range("N17").Value = 20139.0137
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 22108.2106
$ws.Range("S17").Value = 0
$ws.Range("T17").Value = 82101.2093

# Row 18 (Hour 17)
$ws.Range("A18").Value = "'2024-08-15"
$ws.Range("A18").Style = "Normal"
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = 30925
$ws.Range("D18").Value = 6519
$ws.Range("E18").Value = 19370
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 56814
$ws.Range("L18").Value = 30928.0925
$ws.Range("M18").Value = 6519.6519
$ws.Range("N18").Value = 19371.937
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = 0
$ws.Range("S18").Value = 0
$ws.Range("T18").Value = 56819.6814
